$d = $word.ActiveDocument

# Helper: find the range of $searchText starting the search at character
# position $afterPos (searches the story from there to the end of doc).
# The document end is recomputed every call since earlier edits in this
# script change the story's total character count.
function Find-TextRange($searchText, $afterPos) {
    $docEnd = $d.Characters.Count
    $rng = $d.Range($afterPos, $docEnd)
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    return $d.Range($rng.Start, $rng.End)
}

# ---------------------------------------------------------------------
# 1) Title paragraph: drop the leading "Tentative " and advance the
#    version/date text by one day (13/11 -> 1/12).
#    "Tentative Schedule VIP-2021 v. 13/11" -> "Schedule VIP-2021 v. 1/12"
# ---------------------------------------------------------------------
$titleStart = Find-TextRange "Tentative Schedule VIP-2021 v. 13/11" 0

$tStart = $titleStart.Start
$tEnd = $titleStart.End

# Absolute character offsets of each original run within the title range.
$r0 = $tStart             # "Tentative "
$r1 = $tStart + 10        # "Schedule VIP-2"
$r2 = $tStart + 24        # "02"
$r3 = $tStart + 26        # "1"
$r4 = $tStart + 27        # " v. "
$r5 = $tStart + 31        # "13"
$r6 = $tStart + 33        # "/"
$r7 = $tStart + 34        # "1"
$r8 = $tStart + 35        # "1"
$r9 = $tStart + 36        # end

if ($r9 -ne $tEnd) {
    throw "Unexpected title length: expected end $tEnd, computed $r9"
}

# Apply edits right-to-left so untouched offsets stay valid.
$d.Range($r8, $r9).Text = ""      # delete trailing "1"
$d.Range($r7, $r8).Text = "2"     # "1"  -> "2"
$d.Range($r6, $r7).Text = "1"     # "/"  -> "1"
$d.Range($r5, $r6).Text = "/"     # "13" -> "/"
$d.Range($r4, $r5).Text = "1"     # " v. " -> "1"
$d.Range($r3, $r4).Text = " v. "  # "1"  -> " v. "
$d.Range($r2, $r3).Text = "1"     # "02" -> "1"
$d.Range($r1, $r2).Text = "02"    # "Schedule VIP-2" -> "02"
$d.Range($r0, $r1).Text = "Schedule VIP-2"   # "Tentative " -> "Schedule VIP-2"

# ---------------------------------------------------------------------
# 2) "Ex4:  15.12 - 10.1 CBIR- SIO" -> "Ex4:  15.12 - 12.1 CBIR- SIO"
#    (the lone "0" run following "15.12 - 1" becomes "2")
# ---------------------------------------------------------------------
$ex4Anchor = Find-TextRange "15.12" 0
$ex4Ten = Find-TextRange "10.1" $ex4Anchor.End
$d.Range($ex4Ten.Start + 1, $ex4Ten.Start + 2).Text = "2"   # "0" -> "2"

# ---------------------------------------------------------------------
# 3) "Ex5:  10.1 - 19.1  Segmentation- FL"
#       -> "Ex5:  12.1 - 21.1  Segmentation- FL"
#    (the "0" run in the first "10.1" becomes "2"; the "19" run becomes "21")
# ---------------------------------------------------------------------
$ex5Ten = Find-TextRange "10.1" $ex4Ten.End
$ex5Nineteen = Find-TextRange "19.1" $ex5Ten.End

# Right-to-left: the "19" -> "21" edit comes after the "0" -> "2" edit in
# the document, so do it first to keep the earlier offset valid.
$d.Range($ex5Nineteen.Start, $ex5Nineteen.Start + 2).Text = "21"  # "19" -> "21"
$d.Range($ex5Ten.Start + 1, $ex5Ten.Start + 2).Text = "2"         # "0" -> "2"
